# "continua graduação + planejamento"
# The percentage columns (C, E, G) were reformatted to use a comma as the
# decimal separator (pt-BR style) instead of a period, e.g. "55.2%" -> "55,2%".
# Numeric/header columns and row 1 (header) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$percentCols = @("C", "E", "G")

for ($row = 2; $row -le 9; $row++) {
    foreach ($col in $percentCols) {
        $addr = "$col$row"
        $current = $ws.Range($addr).Text
        if ($current -like "*.*") {
            $updated = $current -replace "\.", ","
            $ws.Range($addr).Value = $updated
        }
    }
}
